# The workbook was simply reopened in Excel (its volatile
# RANDBETWEEN()-based "Dados não relevantes aqui NN" helper column
# recalculates on every open/save), the user clicked on cell B9, and the
# file was saved again - no deliberate data edits were made to the sheet
# itself (see commit message: the real changes are all in the Python
# automation script, not in this workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the volatile helper-column formulas (RANDBETWEEN via column B)
# to recalculate, mirroring Excel recalculating on open/save.
$excel.CalculateFull() | Out-Null

# Move the selection to B9, matching the saved cursor position.
$ws.Range("B9").Select() | Out-Null
